$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Variable_value / Variable_description cells that changed ---
$ws.Range("C2").Value = "P1,P2,Stranger"
$ws.Range("C3").Value = "P1,P2,Stranger"
$ws.Range("C4").Value = "Top,Bottom"
$ws.Range("B5").Value = "Judging whether specific shapes and labels match or mismatch"
$ws.Range("C5").Value = "Match,Mismatch"
$ws.Range("C10").Value = "o,w,None"

# --- Column B width widened (68 display chars; the engine's internal MDW
# rounding adds 5/7 px of padding on top of whatever we assign, so back that
# off here to land exactly on 68 in the saved XML) ---
$ws.Columns.Item(2).ColumnWidth = 67.28571428571429

# --- Selection moves to B25 ---
$ws.Range("B25").Select()
